$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# A2 was 1, becomes 0 (the running "TEST #" counter starts at 0 now).
# Downstream A3:A33 are formulas (1+prev) and will recalc automatically.
$ws.Range("A2").Value = 0

# D2/E2 (MIN/MAX for the first test row) become formulas derived from F2
# (NOMINAL), using +/-5% tolerance, instead of hard-coded literals.
$ws.Range("D2").Formula = "=0.95*F2"
$ws.Range("E2").Formula = "=1.05*F2"

# D3:D9 / E3:E9 get the same +/-5% formulas (entered as one range each so
# they form shared-formula groups, matching rows 3 through 9).
$ws.Range("D3:D9").Formula = "=0.95*F3"
$ws.Range("E3:E9").Formula = "=1.05*F3"

# All of D2:F9 picks up the new "0.000" number format (3 decimal places).
$ws.Range("D2:F9").NumberFormat = "0.000"

# Restore the active cell/selection to A3 (was B23).
[void]$ws.Range("A3").Select()
